$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values look numeric (e.g. "6.15"); force those
# specific cells to keep a text format so Excel does not coerce them
# into floating-point numbers (matches the existing inline-string text).
$numericLookingCells = @(
  'D5',
  'D6',
  'D9',
  'D12',
  'D15',
  'D19',
  'D20',
  'D22',
  'D24',
  'D25',
  'D29',
  'D31',
  'D33',
  'D35',
  'D36',
  'D37',
  'D38',
  'D40',
  'D41',
  'D42',
  'D44',
  'D45',
  'D46',
  'D49',
  'D50'
)
foreach ($addr in $numericLookingCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$updates = @(
  @('D2', '57.736.39'),
  @('E2', '  +0.03%  '),
  @('D3', '3.107.70'),
  @('E3', '  +1.43%  '),
  @('E4', '  -0.02%  '),
  @('D5', '524.32'),
  @('E5', '  +0.99%  '),
  @('D6', '141.46'),
  @('E6', '  -0.37%  '),
  @('D8', '3.104.70'),
  @('E8', '  +1.42%  '),
  @('D9', '0.436'),
  @('E9', '  +0.32%  '),
  @('E10', '  -0.28%  '),
  @('E11', '  +1.82%  '),
  @('D12', '0.387'),
  @('E12', '  +3.04%  '),
  @('D13', '3.639.36'),
  @('E13', '  +1.49%  '),
  @('E14', '  +1.53%  '),
  @('D15', '26.27'),
  @('E15', '  +2.34%  '),
  @('E16', '  +0.77%  '),
  @('D17', '57.748.68'),
  @('E17', '  -0.05%  '),
  @('D18', '3.106.26'),
  @('E18', '  +1.87%  '),
  @('D19', '6.15'),
  @('E19', '  +1.15%  '),
  @('D20', '12.83'),
  @('E20', '  -0.01%  '),
  @('E21', '  -0.97%  '),
  @('D22', '337.65'),
  @('E22', '  +2.39%  '),
  @('E23', '  -0.05%  '),
  @('D24', '0.512'),
  @('E24', '  +3.00%  '),
  @('D25', '66.82'),
  @('E25', '  +1.51%  '),
  @('E26', '  -0.28%  '),
  @('E27', '  +0.23%  '),
  @('D28', '0.0₃0924'),
  @('E28', '  +2.59%  '),
  @('D29', '6.54'),
  @('E29', '  +3.06%  '),
  @('E30', '  +0.03%  '),
  @('D31', '7.23'),
  @('E31', '  +0.41%  '),
  @('E32', '  +2.53%  '),
  @('D33', '21.07'),
  @('E33', '  +1.76%  '),
  @('E34', '  +1.36%  '),
  @('D35', '155.43'),
  @('E35', '  +0.56%  '),
  @('D36', '4.66'),
  @('E36', '  +3.69%  '),
  @('D37', '6.12'),
  @('E37', '  +2.96%  '),
  @('D38', '27.24'),
  @('E38', '  -0.73%  '),
  @('E39', '  +2.46%  '),
  @('D40', '0.0666'),
  @('E40', '  -1.27%  '),
  @('B41', 'Filecoin'),
  @('C41', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
  @('D41', '3.94'),
  @('E41', '  +0.70%  '),
  @('B42', 'Stacks'),
  @('C42', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'),
  @('D42', '1.53'),
  @('E42', '  +12.26%  '),
  @('B43', 'RenzoRestakedETH'),
  @('C43', 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'),
  @('D43', '3.145.27'),
  @('E43', '  +1.35%  '),
  @('B44', 'Mantle'),
  @('C44', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
  @('D44', '0.687'),
  @('E44', '  +5.36%  '),
  @('D45', '36.83'),
  @('E45', '  +0.17%  '),
  @('D46', '0.999'),
  @('E46', '  -0.02%  '),
  @('D47', '2.320.48'),
  @('E47', '  +2.62%  '),
  @('E48', '  +1.11%  '),
  @('D49', '0.980'),
  @('E49', '  +6.64%  '),
  @('D50', '20.81'),
  @('E50', '  +0.04%  '),
  @('E51', '  +2.48%  ')
)

foreach ($u in $updates) {
  $ws.Range($u[0]).Value = $u[1]
}
